$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45174 -> 45175) for rows 2 through 13.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
